# Add two new worksheets "Cards" and "Cards1" after the existing "Cube Info"
# sheet, each holding a small card-library table (Spotify Toybox library),
# reusing the header cell style (s="1") already defined in the workbook.

$wb = $excel.ActiveWorkbook
$cubeInfo = $wb.Worksheets.Item(1)

# --- Sheet: Cards -----------------------------------------------------
$cards = $wb.Worksheets.Add($null, $cubeInfo)
$cards.Name = "Cards"

# Copy the existing header style (bold, bordered, centered) from Cube Info!A1
# onto the header row and the numeric A2 cell so a new style isn't created.
$cubeInfo.Range("A1").Copy()
$cards.Range("B1:E1").PasteSpecial(-4122)
$cubeInfo.Range("A1").Copy()
$cards.Range("A2").PasteSpecial(-4122)

$cards.Range("B1").Value = "card_name"
$cards.Range("C1").Value = "card_cid"
$cards.Range("D1").Value = "card_strats"
$cards.Range("E1").Value = "card_tags"

$cards.Range("A2").Value = 0
$cards.Range("B2").Value = "Izzet Charm"
$cards.Range("C2").Value = "UR"
$cards.Range("D2").Value = "Arf, Meow"
$cards.Range("E2").Value = "Commander"

# --- Sheet: Cards1 ------------------------------------------------------
$cards1 = $wb.Worksheets.Add($null, $cards)
$cards1.Name = "Cards1"

$cubeInfo.Range("A1").Copy()
$cards1.Range("B1:E1").PasteSpecial(-4122)
$cubeInfo.Range("A1").Copy()
$cards1.Range("A2").PasteSpecial(-4122)

$cards1.Range("B1").Value = "card_name"
$cards1.Range("C1").Value = "card_cid"
$cards1.Range("D1").Value = "card_strats"
$cards1.Range("E1").Value = "card_tags"

$cards1.Range("A2").Value = 0
$cards1.Range("B2").Value = "Memnite"
$cards1.Range("C2").Value = "UR"
$cards1.Range("D2").Value = "Arf, Meow"
$cards1.Range("E2").Value = "Commander"

$cubeInfo.Select()
